# "bolded new slide title"
# Slide 3's title placeholder ("Title 1") holds the code snippet
#   BankAccount b2 = new SavingsAccount(new Customer("Jane Doe", 2002), 500.0, 5.0);
# split across four runs. Bold the whole title run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Select the entire title text and turn bold on for every run it covers.
$full = $tr.Characters(1, $tr.Length)
$full.Font.Bold = $true
